$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (before the old row 5),
# pushing the existing rows 5-17 down to 7-19.
$ws.Rows("5:6").Insert()

# New row 5: Perejil "Primera" quote dated 44804
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C5").Value = 'Ñuble'
$ws.Range("D5").Value = 44804
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112044
$ws.Range("G5").Value = 'Perejil'
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 750
$ws.Range("L5").Value = 850
$ws.Range("M5").Value = 800
$ws.Range("N5").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O5").Value = 'Región del Maule'
$ws.Range("P5").Value = 800
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 'Hortaliza'

# New row 6: Perejil "Segunda" quote dated 44804
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C6").Value = 'Ñuble'
$ws.Range("D6").Value = 44804
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112044
$ws.Range("G6").Value = 'Perejil'
$ws.Range("H6").Value = 'Sin especificar'
$ws.Range("I6").Value = 'Segunda'
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 650
$ws.Range("L6").Value = 650
$ws.Range("M6").Value = 650
$ws.Range("N6").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O6").Value = 'Región del Maule'
$ws.Range("P6").Value = 650
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 'Hortaliza'
